$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (row 1) to human-readable Spanish text
$ws.Range("A1").Value = "Horas trabajadas"
$ws.Range("B1").Value = "Personas residentes viviendas familiares"
$ws.Range("C1").Value = "Ocupación (1 dígito) descripción"
$ws.Range("D1").Value = "Ocupacion (1 dígito) código"
$ws.Range("E1").Value = "Aragón"

# Fix the data-type row (row 4): the "Horas trabajadas" and
# "Ocupación (1 dígito) descripción" measures are textual (xsd:string),
# not numeric (xsd:int) as previously generated.
$ws.Range("A4").Value = "xsd:string"
$ws.Range("C4").Value = "xsd:string"
